# feat: add 2022-Q4 data
#
# The workbook tracks quarterly holdings. A new quarter ("2022-Q4") is
# prepended: the summary sheet ("总计") gets a new second row for it, and a
# new worksheet named "2022-Q4" (cloned from the "2022-Q3" worksheet, which
# already carries the right layout/formatting) is inserted right after
# "总计" with its own figures. The existing quarter worksheets
# (2022-Q3 / 2022-Q2 / 2021-Q4 / 2021-Q3) keep their names and data - they
# simply shift one tab to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new row 2 for 2022-Q4, pushing the
#    previously-existing rows (2022-Q3 / 2022-Q2 / 2021-Q4 / 2021-Q3) down
#    by one row each (they keep their own values).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# Copy the (now shifted-down) old row 2 formatting onto the fresh row so the
# new row matches the existing style (s="2" on column A, default elsewhere).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)  # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

# Column A is a simple 0-based row index - renumber it now that a row was
# inserted (Rows.Insert shifted the old index values down along with
# everything else, so row 3 is still "0" etc. - fix that up).
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet, cloned from "2022-Q3" (same columns /
#    header / styling), positioned right after "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Stage the new figures as quote-prefixed text in helper cells, then
# paste-special (values only) into the real cells so they land as plain
# text (matching how the source data is stored) without leaving a
# lingering "@" / quote-prefix style behind on the destination cells.
$helper = $q4.Range("J1:M1")
$helper.NumberFormat = "@"
$q4.Range("J1").Value = "0.43"
$q4.Range("K1").Value = "92.90"
$q4.Range("L1").Value = "2.08"
$q4.Range("M1").Value = "0.0089"
$helper.Copy()
$q4.Range("D2").PasteSpecial(-4163)  # xlPasteValues
$helper.Clear()

# ---------------------------------------------------------------------
# 3) Restore the originally-active tab ("2021-Q3") - adding/copying sheets
#    above moved the active selection onto the freshly inserted sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
